$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 0.22510864227085392
$ws.Range("A2").Value = -0.03736982372008768
$ws.Range("A3").Value = -0.0039999999632662764
$ws.Range("A4").Value = -0.007999999934259705
$ws.Range("A5").Value = -0.0029999999628609331
$ws.Range("A6").Value = -0.0041143458827566093
$ws.Range("A7").Value = -0.0099999999101365589
$ws.Range("A8").Value = 0.035311834494863081
$ws.Range("A9").Value = -0.0019999999638553589
$ws.Range("A10").Value = -0.0019999999640063493
$ws.Range("A11").Value = -0.0029999999576348912
$ws.Range("A12").Value = -0.0034999999546161398
$ws.Range("A13").Value = -0.0034999999562304041
$ws.Range("A14").Value = -0.0079999999281659129
$ws.Range("A15").Value = -0.00099999997429467413
$ws.Range("A16").Value = 0.029032666440026222
$ws.Range("A17").Value = -0.0019999999688273817
$ws.Range("A18").Value = -0.0039999999557558397
$ws.Range("A19").Value = -0.0039999999727582392
$ws.Range("A20").Value = -0.0039999999707269751
$ws.Range("A21").Value = -0.0039999999704223299
$ws.Range("A22").Value = -0.0039999999701594291
$ws.Range("A23").Value = -0.0049999999530614403
$ws.Range("A24").Value = -0.019999999849083849
$ws.Range("A25").Value = -0.051675097983191698
$ws.Range("A26").Value = -0.0024999999555728181
$ws.Range("A27").Value = -0.0024999999527568484
$ws.Range("A28").Value = -0.0019999999431092874
$ws.Range("A29").Value = -0.0069999999015291081
$ws.Range("A30").Value = -0.059999999557388328
$ws.Range("A31").Value = -0.0069999998937451124
$ws.Range("A32").Value = -0.0099999998740347706
$ws.Range("A33").Value = -0.0039999999123025987
